$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Sander")

# Insert a new column before column B (shifts existing B.. columns to the right)
$ws.Columns.Item(2).Insert()

# New column inherits column A's width
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Header
$ws.Cells.Item(1, 2).Value = "Bedrijf2"

# Anonymized company labels A, B, C, ... Z, AA, AB, AC, AD for rows 2..31
$labels = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $labels[$i]
}

# Match cursor position left after the edit
$ws.Range("B32").Select()
